$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is shrinking from a 4x3 grid (A1:D3) to a 3x2 grid (A1:C2):
# drop the now-unused column D and the third data row.
$ws.Range("D1:D3").EntireColumn.Delete()
$ws.Range("A3:C3").EntireRow.Delete()

# Row 1 holds text labels (stored as strings, not numbers) for the two
# remaining columns that change. Use a leading apostrophe so the
# numeric-looking strings "-1"/"-2" are kept as text instead of being
# reinterpreted as numbers. C1's label ("0") is unchanged, so it is left
# untouched.
$ws.Range("A1").Formula = "'-1"
$ws.Range("B1").Formula = "'-2"

# Forcing text on A1:B1 allocates a new cell style, which would otherwise
# diverge from the sheet's original formatting. Restore the original
# border/font/alignment by copying the (untouched) format from C1, which
# still carries the original style.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 holds the numeric data values.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2.5
